$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Force the cell to remain a text value even when $value looks numeric
    # (e.g. "218.05"), matching the source inlineStr cells, then restore the
    # default (unstyled) cell style so no spurious style index is introduced.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "26.130.95"
$ws.Range("E2").Value = "  -0.62%  "
Set-TextCell $ws.Range("D3") "1.655.85"
$ws.Range("E3").Value = "  -0.77%  "
Set-TextCell $ws.Range("D5") "218.05"
$ws.Range("E5").Value = "  -0.26%  "
Set-TextCell $ws.Range("D6") "0.5290"
$ws.Range("E6").Value = "  +0.81%  "
Set-TextCell $ws.Range("D7") "1.003"
$ws.Range("E7").Value = "  -0.30%  "
Set-TextCell $ws.Range("D8") "0.2606"
$ws.Range("E8").Value = "  -2.48%  "
Set-TextCell $ws.Range("D9") "0.06337"
$ws.Range("E9").Value = "  +0.06%  "
Set-TextCell $ws.Range("D10") "20.41"
$ws.Range("E10").Value = "  -2.99%  "
Set-TextCell $ws.Range("D11") "0.07783"
$ws.Range("E11").Value = "  +0.28%  "
Set-TextCell $ws.Range("D12") "4.496"
$ws.Range("E12").Value = "  +0.99%  "
Set-TextCell $ws.Range("D13") "1.642.03"
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("E14").Value = "  +0.10%  "
Set-TextCell $ws.Range("D15") "0.0₅8164"
$ws.Range("E15").Value = "  -1.40%  "
Set-TextCell $ws.Range("D16") "65.46"
$ws.Range("E16").Value = "  +0.63%  "
Set-TextCell $ws.Range("D17") "26.135.14"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("E18").Value = "  -0.37%  "
Set-TextCell $ws.Range("D19") "4.563"
$ws.Range("E19").Value = "  -2.51%  "
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("E21").Value = "  -0.72%  "
Set-TextCell $ws.Range("D22") "6.024"
$ws.Range("E22").Value = "  -0.94%  "
Set-TextCell $ws.Range("D23") "1.003"
$ws.Range("E23").Value = "  -0.45%  "
Set-TextCell $ws.Range("D24") "141.81"
$ws.Range("E24").Value = "  +1.29%  "
Set-TextCell $ws.Range("D25") "0.1252"
$ws.Range("E25").Value = "  +0.98%  "
Set-TextCell $ws.Range("D26") "7.266"
$ws.Range("E26").Value = "  +0.77%  "
Set-TextCell $ws.Range("D27") "16.19"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  +1.70%  "
Set-TextCell $ws.Range("D29") "0.05935"
$ws.Range("E29").Value = "  -3.84%  "
Set-TextCell $ws.Range("D30") "1.279"
$ws.Range("E30").Value = "  -0.31%  "
Set-TextCell $ws.Range("D31") "3.513"
$ws.Range("E31").Value = "  -2.22%  "
Set-TextCell $ws.Range("D32") "3.248"
$ws.Range("E32").Value = "  -1.57%  "
Set-TextCell $ws.Range("D33") "1.574"
$ws.Range("E33").Value = "  -3.60%  "
$ws.Range("E34").Value = "  +0.15%  "
Set-TextCell $ws.Range("D35") "0.9492"
$ws.Range("E35").Value = "  -2.52%  "
Set-TextCell $ws.Range("D36") "2.410"
$ws.Range("E36").Value = "  -0.69%  "
Set-TextCell $ws.Range("D37") "0.5668"
$ws.Range("E37").Value = "  -1.68%  "
Set-TextCell $ws.Range("D38") "0.01609"
$ws.Range("E38").Value = "  -0.06%  "
Set-TextCell $ws.Range("D39") "5.815"
$ws.Range("E39").Value = "  -3.49%  "
Set-TextCell $ws.Range("D40") "0.8475"
$ws.Range("E40").Value = "  -1.44%  "
Set-TextCell $ws.Range("D41") "1.003"
$ws.Range("E41").Value = "  -0.22%  "
Set-TextCell $ws.Range("D42") "102.46"
$ws.Range("E42").Value = "  +2.11%  "
Set-TextCell $ws.Range("D43") "1.021.81"
$ws.Range("E43").Value = "  -0.24%  "
Set-TextCell $ws.Range("D44") "1.800.01"
$ws.Range("E44").Value = "  -0.63%  "
Set-TextCell $ws.Range("D45") "57.20"
$ws.Range("E45").Value = "  -0.94%  "
Set-TextCell $ws.Range("D46") "1.008"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("E47").Value = "  +1.61%  "
$ws.Range("E48").Value = "  -0.77%  "
Set-TextCell $ws.Range("D49") "0.05150"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("E50").Value = "  -3.78%  "
Set-TextCell $ws.Range("D51") "0.09692"
$ws.Range("E51").Value = "  -1.09%  "
